# Iteration 2 document build-out.
#
# Strategy:
#  - Use Selection.TypeText / TypeParagraph to add all the paragraphs in
#    document order, starting right after the existing bookmarked
#    paragraph.
#  - For paragraphs that must end up fully empty (serialized as a bare
#    <w:p/> with no run at all) we type a single placeholder character,
#    then shrink that paragraph's Range by one unit (dropping the
#    paragraph mark) and clear the Range's Text. That reliably collapses
#    the paragraph back down to a true empty <w:p/> instead of leaving a
#    stray empty <w:r/>.
#  - For the "Step 3" paragraph, which needs <w:proofErr> grammar-check
#    markers wrapped around the word "One", we build the exact paragraph
#    OOXML by hand and drop it in with Range.InsertXML on a collapsed
#    range inside a placeholder paragraph (InsertXML on a collapsed
#    end-of-paragraph range replaces that paragraph's contents).

$d   = $word.ActiveDocument
$sel = $word.Selection

function Clear-ParagraphToEmpty($paraIndex) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.End = $r.End - 1
    $r.Text = ""
}

# Start at the very end of the document (inside the pre-existing
# bookmarked paragraph).
$sel.EndKey(6)

# Paragraph 1 (existing paragraph with the _GoBack bookmark) gains the
# title run.
$sel.TypeText("Iteration 2: Identifying Structures to Support Primary Functionality")

# Paragraph 2: empty.
$sel.TypeParagraph()
$sel.TypeText(".")

# Paragraph 3.
$sel.TypeParagraph()
$sel.TypeText("The goal of this iteration is to address the general architectural concern of identifying structures that support the primary functionality of the application.")

# Paragraph 4: empty.
$sel.TypeParagraph()
$sel.TypeText(".")

# Paragraph 5.
$sel.TypeParagraph()
$sel.TypeText("Step 2: Establish Iteration Goal by Selecting Drivers")

# Paragraph 6.
$sel.TypeParagraph()
$sel.TypeText("The following primary use cases will be addressed:")

# Paragraph 7.
$sel.TypeParagraph()
$sel.TypeText("UC-2: Create/remove content")

# Paragraph 8.
$sel.TypeParagraph()
$sel.TypeText("UC-4: View Course Information")

# Paragraph 9: empty.
$sel.TypeParagraph()
$sel.TypeText(".")

# Paragraph 10: placeholder, content replaced below via InsertXML so we
# can embed the <w:proofErr> markers around "One".
$sel.TypeParagraph()
$sel.TypeText(".")

# Paragraph 11.
$sel.TypeParagraph()
$sel.TypeText(" In this iteration we will be refining the server" + [char]0x2019 + "s architecture in the Rich Client Application Architecture. ")

# Paragraphs 12-16: empty.
$sel.TypeParagraph()
$sel.TypeText(".")
$sel.TypeParagraph()
$sel.TypeText(".")
$sel.TypeParagraph()
$sel.TypeText(".")
$sel.TypeParagraph()
$sel.TypeText(".")
$sel.TypeParagraph()
$sel.TypeText(".")

# Now clean up every placeholder "." paragraph back down to a true empty
# <w:p/>. Do this from the bottom of the document upward so indices
# earlier in the document are not disturbed by the edits.
Clear-ParagraphToEmpty 16
Clear-ParagraphToEmpty 15
Clear-ParagraphToEmpty 14
Clear-ParagraphToEmpty 13
Clear-ParagraphToEmpty 12
Clear-ParagraphToEmpty 9
Clear-ParagraphToEmpty 4
Clear-ParagraphToEmpty 2

# Replace paragraph 10's placeholder with the real content, including the
# <w:proofErr> grammar-check markers Word leaves around "One".
$p10 = $d.Paragraphs(10)
$r10 = $p10.Range
$r10.Collapse(0)
$xml10 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:t xml:space="preserve">Step 3: Choose </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>One</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> or More elements of the system to refine</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r10.InsertXML($xml10)
